$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.536.55"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").Value = "1.755.03"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4585"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.30%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3554"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.89%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07456"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.48"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.086"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.56%  "
$ws.Range("E12").Value = "  +0.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.78"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.002"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.146"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.40%  "
$ws.Range("D16").Value = "1.753.04"
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.49"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.65%  "
$ws.Range("E18").Value = "  -1.10%  "
$ws.Range("E19").Value = "  +0.23%  "
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.83%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.735"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.11%  "
$ws.Range("D23").Value = "27.587.20"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.068"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.92%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.81"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.43%  "
$ws.Range("E27").Value = "  -1.21%  "
$ws.Range("D28").Value = "1.954.28"
$ws.Range("E28").Value = "  +0.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.128"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.081"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09218"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.22%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.662"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.516"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.57%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.75"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.33%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02278"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2091"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06014"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.45%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6270"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.72%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.922"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.20%  "
$ws.Range("E41").Value = "  -1.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.387"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.768"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.51%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.20"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.80%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.716"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5861"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "121.98"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.935"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06897"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.127"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.09"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.60%  "
